# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Wed Sep 20 19:23:34 UTC 2023 with GitHub Actions"
#
# Prices/volumes moved and a few coins swapped table position (rows 15/16,
# 40/41, 43/44). Numeric-looking price strings are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's inlineStr
# cells) instead of silently coercing them to numbers and dropping
# significant trailing zeros (e.g. "215.00" -> 215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.059.80"
$ws.Range("E2").Value = "  -0.38%  "

# Row 3
$ws.Range("D3").Value = "1.623.86"
$ws.Range("E3").Value = "  -1.33%  "

# Row 4
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").Value = "'215.00"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  -0.78%  "

# Row 9
$ws.Range("D9").Value = "'0.0623"
$ws.Range("E9").Value = "  -0.58%  "

# Row 10
$ws.Range("D10").Value = "'20.09"
$ws.Range("E10").Value = "  +0.88%  "

# Row 11
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  -0.12%  "

# Row 12
$ws.Range("D12").Value = "1.633.03"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("D13").Value = "'4.12"
$ws.Range("E13").Value = "  -0.46%  "

# Row 14
$ws.Range("D14").Value = "'0.541"
$ws.Range("E14").Value = "  +0.03%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "27.120.73"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'64.32"
$ws.Range("E16").Value = "  -4.78%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0738"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "'215.93"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").Value = "'6.89"
$ws.Range("E20").Value = "  +1.00%  "

# Row 21
$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22
$ws.Range("E22").Value = "  -6.14%  "

# Row 23
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  -2.27%  "

# Row 24
$ws.Range("D24").Value = "'147.63"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26
$ws.Range("D26").Value = "'7.26"
$ws.Range("E26").Value = "  -3.94%  "

# Row 27
$ws.Range("D27").Value = "'0.118"
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("D28").Value = "'15.58"
$ws.Range("E28").Value = "  -1.01%  "

# Row 29
$ws.Range("D29").Value = "'0.0504"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("E30").Value = "  -0.64%  "

# Row 31
$ws.Range("D31").Value = "'3.36"
$ws.Range("E31").Value = "  -0.99%  "

# Row 32
$ws.Range("D32").Value = "'2.98"
$ws.Range("E32").Value = "  -1.55%  "

# Row 33
$ws.Range("D33").Value = "1.334.48"
$ws.Range("E33").Value = "  +5.69%  "

# Row 34
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  -0.59%  "

# Row 35
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("D36").Value = "'0.0175"
$ws.Range("E36").Value = "  -1.22%  "

# Row 37
$ws.Range("D37").Value = "'0.541"
$ws.Range("E37").Value = "  -1.15%  "

# Row 38
$ws.Range("D38").Value = "'0.845"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.804"
$ws.Range("E40").Value = "  -0.44%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.24"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
$ws.Range("D42").Value = "'64.22"
$ws.Range("E42").Value = "  +3.79%  "

# Row 43
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.770.52"
$ws.Range("E43").Value = "  -0.96%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.24"
$ws.Range("E44").Value = "  -3.47%  "

# Row 45
$ws.Range("D45").Value = "'90.44"
$ws.Range("E45").Value = "  -1.11%  "

# Row 46
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").Value = "'0.821"
$ws.Range("E47").Value = "  +22.44%  "

# Row 48
$ws.Range("D48").Value = "0.0₆0100"
$ws.Range("E48").Value = "  -6.62%  "

# Row 49
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("D50").Value = "'0.0986"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51
$ws.Range("D51").Value = "'7.53"
$ws.Range("E51").Value = "  -1.20%  "
